# "update to adapt new batch of data"
# Add two new mapping rows (study.study_access, file.file_access) to the
# "Must have properties" sheet, right below the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Must have properties")

# Row 17: study.study_access -> study / study_access
$ws.Cells.Item(17, 1).Value = "study.study_access"
$ws.Cells.Item(17, 2).Value = "study"
$ws.Cells.Item(17, 3).Value = "study_access"

# Row 18: file.file_access -> file / file_access
$ws.Cells.Item(18, 1).Value = "file.file_access"
$ws.Cells.Item(18, 2).Value = "file"
$ws.Cells.Item(18, 3).Value = "file_access"

# Match the author's last selection/cursor position on save.
$ws.Range("C21").Select() | Out-Null
